$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths: A and B get wider, and a new column H is introduced.
# (ColumnWidth is quantized to whole pixels internally, so we dial in the
# value that lands closest to the target stored width.)
# ---------------------------------------------------------------------------
$ws.Range("A1").EntireColumn.ColumnWidth = 29.166666666666668
$ws.Range("B1").EntireColumn.ColumnWidth = 35.0
$ws.Range("H1").EntireColumn.ColumnWidth = 18.5

# ---------------------------------------------------------------------------
# New header cell H1 = "confirmPassword", formatted like the other header
# cells (copy format from G1 first, so it picks up the shared header style).
# ---------------------------------------------------------------------------
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "confirmPassword"

# ---------------------------------------------------------------------------
# H2 = "RANDOM_PASSWORD" (same text as G2), formatted like G2.
# ---------------------------------------------------------------------------
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Value = "RANDOM_PASSWORD"

# ---------------------------------------------------------------------------
# New row 3: RegisterUser_EmailAlreadyRegistered test case.
# Plain text columns (A3, C3, D3, E3, G3, H3) reuse the "normal" data style
# (copied from A2/C2/etc.), while B3 keeps the apostrophe/quote-prefixed
# style from B2, and F3 needs to store a numeric-looking phone number as
# literal text but WITHOUT the quote-prefix flag (it uses the plain style).
# ---------------------------------------------------------------------------

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Value = "RegisterUser_EmailAlreadyRegistered"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Value = "'Warning: E-Mail Address is already registered!"

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Value = "Jane"

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Value = "Doe"

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Value = "abcd@gmail.com"

# F3 needs text type ("1234567891") under the plain style (no quote-prefix).
# Build the text value in a scratch cell that already carries the
# quote-prefix style (so the apostrophe-forced text reuses that style
# cleanly), then copy only the *value* across onto F3 (which already has
# the plain style applied), and finally wipe the scratch cell.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Value = "'1234567891"

$ws.Range("C2").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null

$ws.Range("Z1").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4163) | Out-Null

$ws.Range("Z1").Clear() | Out-Null

$ws.Range("G2").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").Value = "RANDOM_PASSWORD"

$ws.Range("G2").Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4122) | Out-Null
$ws.Range("H3").Value = "RANDOM_PASSWORD"
